$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "id" column (column A). This shifts every other column
#    one place to the left (B->A, C->B, ... H->G) and drops the now
#    unused "id" shared string automatically on save.
$ws.Columns.Item(1).Delete()

# 2. Rename the header row to the new, human readable column titles.
$ws.Cells.Item(1, 1).Value = "First name"
$ws.Cells.Item(1, 2).Value = "Last name"
$ws.Cells.Item(1, 3).Value = "Phone number"
$ws.Cells.Item(1, 4).Value = "Address"
$ws.Cells.Item(1, 5).Value = "Passport details"
$ws.Cells.Item(1, 6).Value = "Date of Birth"
$ws.Cells.Item(1, 7).Value = "Discount code"

# 3. Reformat the "Date of Birth" column (now column F) from the verbose
#    JS Date.toString() style strings into "DD : MM : YYYY" strings.
$ws.Cells.Item(2, 6).Value = "18 : 01 : 1990"
$ws.Cells.Item(3, 6).Value = "24 : 02 : 1980"
$ws.Cells.Item(4, 6).Value = "20 : 03 : 1996"
$ws.Cells.Item(5, 6).Value = "05 : 05 : 1978"
$ws.Cells.Item(6, 6).Value = "03 : 06 : 1988"
$ws.Cells.Item(7, 6).Value = "13 : 07 : 1967"
$ws.Cells.Item(8, 6).Value = "04 : 08 : 1987"
$ws.Cells.Item(9, 6).Value = "28 : 09 : 1976"
$ws.Cells.Item(10, 6).Value = "01 : 10 : 1993"
$ws.Cells.Item(11, 6).Value = "15 : 11 : 1995"

# 4. Add a brand new data row (row 12) with test/demo data entered via
#    the admin panel.
$ws.Cells.Item(12, 1).Value = "цу"
$ws.Cells.Item(12, 2).Value = "ук"
$ws.Cells.Item(12, 3).Value = "ку"
$ws.Cells.Item(12, 4).Value = "цуацуацуацауцацуацуа цуацкккккуацуауууууууу"
$ws.Cells.Item(12, 5).Value = "23вуца"
$ws.Cells.Item(12, 6).Value = "14 : 06 : 2019"
$ws.Cells.Item(12, 7).Value = "21квцу"
